$d = $word.ActiveDocument

# Edit 1: Software Engineering skills paragraph - rewrite the tail of the
# sentence about design patterns / REST APIs / HTML / CSS, and remove the
# duplicated sentence that repeated "Have also worked with REST APIs...".
$find1 = ", as well as various design patterns associated with object-oriented programming. Have also worked with REST APIs, HTML, and CSS, and can quickly learn other tools and skills as necessary."
$replace1 = ", and have also worked with REST APIs, HTML, and CSS. Can quickly learn other tools and skills as necessary, and have been doing so alongside my job search."
$found1 = $d.Content.Find.Execute($find1, $false, $true, $false, $false, $false, $true, 1, $false, $replace1, 2)
if (-not $found1) { throw "Edit 1 text not found" }

# Edit 2: Personal projects bullet - move the "personal website" mention
# earlier in the sentence and drop the now-redundant trailing sentence.
$find2 = "Worked on several personal projects while searching for a job and built various hard and soft skills to help with this. My most recent project is the personal website linked at the top of this resume."
$replace2 = "Worked on several projects such as the personal website linked at the top of this resume while searching for a job and built various hard and soft skills to help with this."
$found2 = $d.Content.Find.Execute($find2, $false, $true, $false, $false, $false, $true, 1, $false, $replace2, 2)
if (-not $found2) { throw "Edit 2 text not found" }

# Edit 3: the document is shorter now, so the stale cached
# lastRenderedPageBreak marker on the final "Dean's List" paragraph no
# longer applies. Touching that run via Find/Replace (text unchanged)
# makes Word regenerate it without the obsolete marker.
$find3 = "Made the Dean’s List every semester except for my first two at Marshall"
$found3 = $d.Content.Find.Execute($find3, $false, $true, $false, $false, $false, $true, 1, $false, $find3, 2)
if (-not $found3) { throw "Edit 3 text not found" }
